$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple single-run cell text replacements (by row index, before any row add/delete) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "44"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00047"
$t.Rows.Item(7).Cells.Item(1).Range.Text = "0.00018"

# --- Delete the row that held "0.00002" (row 8) ---
$t.Rows.Item(8).Delete()

# After deletion, rows shift up by one: the row that held "0.00007" is now row 9,
# "0.00009" is row 10, "0.00084" is row 11, "100.0" is row 12.
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00023"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00025"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00034"

# --- Insert a brand-new row (with "0.00796") right after row 11, before row 12 ("100.0") ---
$refRow = $t.Rows.Item(12)
$newRow = $t.Rows.Add($refRow)
$newRow.Cells.Item(1).Range.Text = "0.00796"

# --- Collapse the final three multi-run rows into single-run cells ---
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "99.94"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "12"
